$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly (closest achievable to target OOXML width 23.7109375)
$ws.Range("C1").ColumnWidth = 22.8

# Update data rows 2-22 for columns A, B, C, H, I, J per forecast refresh
$ws.Cells.Item(2, 1).Value = 255
$ws.Cells.Item(2, 2).Value = 9210
$ws.Cells.Item(2, 3).Value = 30.70000076293945
$ws.Cells.Item(2, 8).Value = -23
$ws.Cells.Item(2, 9).Value = -56647
$ws.Cells.Item(2, 10).Value = 43

$ws.Cells.Item(3, 1).Value = 62
$ws.Cells.Item(3, 2).Value = 646
$ws.Cells.Item(3, 3).Value = 14.04347801208496
$ws.Cells.Item(3, 8).Value = -41
$ws.Cells.Item(3, 9).Value = -89845
$ws.Cells.Item(3, 10).Value = 145

$ws.Cells.Item(4, 1).Value = 53
$ws.Cells.Item(4, 2).Value = 898
$ws.Cells.Item(4, 3).Value = 7.126984119415283
$ws.Cells.Item(4, 8).Value = 35
$ws.Cells.Item(4, 9).Value = 31251
$ws.Cells.Item(4, 10).Value = 2

$ws.Cells.Item(5, 1).Value = 166
$ws.Cells.Item(5, 2).Value = 9564
$ws.Cells.Item(5, 3).Value = 13.28333377838135
$ws.Cells.Item(5, 8).Value = 63
$ws.Cells.Item(5, 9).Value = 19178
$ws.Cells.Item(5, 10).Value = 313

$ws.Cells.Item(6, 1).Value = 288
$ws.Cells.Item(6, 2).Value = 34271
$ws.Cells.Item(6, 3).Value = 35.73618316650391
$ws.Cells.Item(6, 8).Value = 158
$ws.Cells.Item(6, 9).Value = 136886
$ws.Cells.Item(6, 10).Value = 328

$ws.Cells.Item(7, 1).Value = 217
$ws.Cells.Item(7, 2).Value = 58652
$ws.Cells.Item(7, 3).Value = 64.3114013671875
$ws.Cells.Item(7, 8).Value = 199
$ws.Cells.Item(7, 9).Value = 196804
$ws.Cells.Item(7, 10).Value = 375

$ws.Cells.Item(8, 1).Value = 131
$ws.Cells.Item(8, 2).Value = 10614
$ws.Cells.Item(8, 3).Value = 23.43046379089355
$ws.Cells.Item(8, 8).Value = 54
$ws.Cells.Item(8, 9).Value = 6904
$ws.Cells.Item(8, 10).Value = 212

$ws.Cells.Item(9, 1).Value = 142
$ws.Cells.Item(9, 2).Value = 26795
$ws.Cells.Item(9, 3).Value = 79.74702453613281
$ws.Cells.Item(9, 8).Value = 21
$ws.Cells.Item(9, 9).Value = -25653
$ws.Cells.Item(9, 10).Value = 172

$ws.Cells.Item(10, 1).Value = 99
$ws.Cells.Item(10, 2).Value = 5100
$ws.Cells.Item(10, 3).Value = 21.51898765563965
$ws.Cells.Item(10, 8).Value = 30
$ws.Cells.Item(10, 9).Value = -17228
$ws.Cells.Item(10, 10).Value = 152

$ws.Cells.Item(11, 1).Value = 131
$ws.Cells.Item(11, 2).Value = 8732
$ws.Cells.Item(11, 3).Value = 16.50661659240723
$ws.Cells.Item(11, 8).Value = 62
$ws.Cells.Item(11, 9).Value = 10476
$ws.Cells.Item(11, 10).Value = 185

$ws.Cells.Item(12, 1).Value = 260
$ws.Cells.Item(12, 2).Value = 61888
$ws.Cells.Item(12, 3).Value = 149.1277160644531
$ws.Cells.Item(12, 8).Value = 106
$ws.Cells.Item(12, 9).Value = 47479
$ws.Cells.Item(12, 10).Value = 179

$ws.Cells.Item(13, 1).Value = 228
$ws.Cells.Item(13, 2).Value = 93465
$ws.Cells.Item(13, 3).Value = 113.2909088134766
$ws.Cells.Item(13, 8).Value = 190
$ws.Cells.Item(13, 9).Value = 148930
$ws.Cells.Item(13, 10).Value = 290

$ws.Cells.Item(14, 1).Value = 1067
$ws.Cells.Item(14, 2).Value = 786198
$ws.Cells.Item(14, 3).Value = 636.5975952148438
$ws.Cells.Item(14, 8).Value = 573
$ws.Cells.Item(14, 9).Value = 686201
$ws.Cells.Item(14, 10).Value = 543

$ws.Cells.Item(15, 1).Value = 387
$ws.Cells.Item(15, 2).Value = 134647
$ws.Cells.Item(15, 3).Value = 260.4390563964844
$ws.Cells.Item(15, 8).Value = 76
$ws.Cells.Item(15, 9).Value = 29505
$ws.Cells.Item(15, 10).Value = 230

$ws.Cells.Item(16, 1).Value = 64
$ws.Cells.Item(16, 2).Value = 1775
$ws.Cells.Item(16, 3).Value = 18.88297843933105
$ws.Cells.Item(16, 8).Value = 53
$ws.Cells.Item(16, 9).Value = 66922
$ws.Cells.Item(16, 10).Value = -31

$ws.Cells.Item(17, 1).Value = 85
$ws.Cells.Item(17, 2).Value = 2338
$ws.Cells.Item(17, 3).Value = 19.1639347076416
$ws.Cells.Item(17, 8).Value = 43
$ws.Cells.Item(17, 9).Value = 39391
$ws.Cells.Item(17, 10).Value = -31

$ws.Cells.Item(18, 1).Value = 118
$ws.Cells.Item(18, 2).Value = 18388
$ws.Cells.Item(18, 3).Value = 75.36065673828125
$ws.Cells.Item(18, 8).Value = 88
$ws.Cells.Item(18, 9).Value = 39299
$ws.Cells.Item(18, 10).Value = 334

$ws.Cells.Item(19, 1).Value = 259
$ws.Cells.Item(19, 2).Value = 81724
$ws.Cells.Item(19, 3).Value = 213.3785858154297
$ws.Cells.Item(19, 8).Value = 104
$ws.Cells.Item(19, 9).Value = 64934
$ws.Cells.Item(19, 10).Value = 268

$ws.Cells.Item(20, 1).Value = 44
$ws.Cells.Item(20, 2).Value = 12674
$ws.Cells.Item(20, 3).Value = 22.7132625579834
$ws.Cells.Item(20, 8).Value = 46
$ws.Cells.Item(20, 9).Value = 1597
$ws.Cells.Item(20, 10).Value = 105

$ws.Cells.Item(21, 1).Value = 57
$ws.Cells.Item(21, 2).Value = 9198
$ws.Cells.Item(21, 3).Value = 16.9392261505127
$ws.Cells.Item(21, 8).Value = 60
$ws.Cells.Item(21, 9).Value = 17716
$ws.Cells.Item(21, 10).Value = 244

$ws.Cells.Item(22, 1).Value = 339
$ws.Cells.Item(22, 2).Value = 185155
$ws.Cells.Item(22, 3).Value = 261.9937133789062
$ws.Cells.Item(22, 8).Value = 198
$ws.Cells.Item(22, 9).Value = 197834
$ws.Cells.Item(22, 10).Value = 395
